# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a few cell clears) to match target diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 44996.668  # H87
$ws.Cells.Item(87, 10).Value = 44996.668  # J87
$ws.Cells.Item(87, 12).Value = 44996.668  # L87
$ws.Cells.Item(87, 14).Value = -47492.668  # N87
$ws.Cells.Item(90, 8).Value = 44996.668  # H90
$ws.Cells.Item(90, 10).Value = 44996.668  # J90
$ws.Cells.Item(90, 12).Value = 134990.004  # L90
$ws.Cells.Item(90, 14).Value = -147470.004  # N90
$ws.Cells.Item(115, 8).Value = 125  # H115
$ws.Cells.Item(115, 9).Value = 125  # I115
$ws.Cells.Item(115, 11).Value = 375  # K115
$ws.Cells.Item(115, 13).Value = 1192  # M115
$ws.Cells.Item(138, 8).Value = 2632.889  # H138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 29361.5  # H2
$ws.Cells.Item(2, 9).Value = 2499  # I2
$ws.Cells.Item(2, 11).Value = 2499  # K2
$ws.Cells.Item(2, 13).Value = -2386  # M2
$ws.Cells.Item(5, 8).Value = 139.8  # H5
$ws.Cells.Item(5, 9).Value = 139.8  # I5
$ws.Cells.Item(5, 11).Value = 139.8  # K5
$ws.Cells.Item(5, 13).Value = -27.80000000000001  # M5
$ws.Cells.Item(45, 8).Value = 931.8570999999999  # H45
$ws.Cells.Item(45, 10).Value = 1213.3334  # J45
$ws.Cells.Item(45, 12).Value = 1213.3334  # L45
$ws.Cells.Item(45, 14).Value = -1967.3334  # N45
$ws.Cells.Item(74, 8).Value = 3302.2  # H74
$ws.Cells.Item(74, 9).Value = 3049.25  # I74
$ws.Cells.Item(74, 11).Value = 3049.25  # K74
$ws.Cells.Item(74, 13).Value = -2175.25  # M74
$ws.Cells.Item(77, 8).Value = 3302.2  # H77
$ws.Cells.Item(77, 9).Value = 3049.25  # I77
$ws.Cells.Item(77, 11).Value = 15246.25  # K77
$ws.Cells.Item(77, 13).Value = -10878.25  # M77
$ws.Cells.Item(102, 8).Value = 1929.1666  # H102
$ws.Cells.Item(102, 9).Value = 1929.1666  # I102
$ws.Cells.Item(102, 11).Value = 1929.1666  # K102
$ws.Cells.Item(102, 13).Value = -307.1666  # M102
$ws.Cells.Item(110, 8).Value = 2629.3333  # H110
$ws.Cells.Item(110, 9).Value = 2174.4  # I110
$ws.Cells.Item(110, 10).Value = 3198  # J110
$ws.Cells.Item(110, 11).Value = 2174.4  # K110
$ws.Cells.Item(110, 12).Value = 3198  # L110
$ws.Cells.Item(110, 13).Value = -129.4000000000001  # M110
$ws.Cells.Item(110, 14).Value = -7288  # N110
$ws.Cells.Item(116, 8).Value = 29361.5  # H116
$ws.Cells.Item(116, 9).Value = 2499  # I116
$ws.Cells.Item(116, 11).Value = 2499  # K116
$ws.Cells.Item(116, 13).Value = -205  # M116
$ws.Cells.Item(122, 8).Value = 2202.1538  # H122
$ws.Cells.Item(122, 9).Value = 1787.75  # I122
$ws.Cells.Item(122, 10).Value = 2865.2  # J122
$ws.Cells.Item(122, 11).Value = 5363.25  # K122
$ws.Cells.Item(122, 12).Value = 8595.599999999999  # L122
$ws.Cells.Item(122, 13).Value = -2913.25  # M122
$ws.Cells.Item(122, 14).Value = -13495.6  # N122
$ws.Cells.Item(132, 8).Value = 2739.7  # H132
$ws.Cells.Item(132, 9).Value = 2739.7  # I132
$ws.Cells.Item(132, 11).Value = 8219.099999999999  # K132
$ws.Cells.Item(132, 13).Value = -5689.099999999999  # M132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 29361.5  # H3
$ws.Cells.Item(3, 9).Value = 2499  # I3
$ws.Cells.Item(3, 11).Value = 2499  # K3
$ws.Cells.Item(3, 13).Value = -2385  # M3
$ws.Cells.Item(4, 8).Value = 139.8  # H4
$ws.Cells.Item(4, 9).Value = 139.8  # I4
$ws.Cells.Item(4, 11).Value = 139.8  # K4
$ws.Cells.Item(4, 13).Value = -24.80000000000001  # M4
$ws.Cells.Item(20, 8).Value = 1587  # H20
$ws.Cells.Item(20, 9).Value = 0  # I20
$ws.Cells.Item(20, 11).Value = 0  # K20
$ws.Cells.Item(22, 8).Value = 283.33334  # H22
$ws.Cells.Item(22, 9).Value = 283.33334  # I22
$ws.Cells.Item(22, 11).Value = 283.33334  # K22
$ws.Cells.Item(22, 13).Value = -110.33334  # M22
$ws.Cells.Item(94, 8).Value = 1999.6666  # H94
$ws.Cells.Item(94, 9).Value = 1999.6666  # I94
$ws.Cells.Item(94, 10).Value = 0  # J94
$ws.Cells.Item(94, 11).Value = 1999.6666  # K94
$ws.Cells.Item(94, 12).Value = 0  # L94
$ws.Cells.Item(94, 13).Value = -1548.6666  # M94
$ws.Cells.Item(105, 8).Value = 2196.8667  # H105
$ws.Cells.Item(105, 9).Value = 2189.5  # I105
$ws.Cells.Item(105, 11).Value = 2189.5  # K105
$ws.Cells.Item(105, 13).Value = -442.5  # M105
$ws.Cells.Item(107, 8).Value = 762.4  # H107
$ws.Cells.Item(107, 9).Value = 762.4  # I107
$ws.Cells.Item(107, 11).Value = 762.4  # K107
$ws.Cells.Item(107, 13).Value = 1157.6  # M107
$ws.Cells.Item(20, 13).ClearContents()  # M20
$ws.Cells.Item(94, 14).ClearContents()  # N94

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1590.6428  # H31
$ws.Cells.Item(31, 9).Value = 1590.6428  # I31
$ws.Cells.Item(31, 11).Value = 1590.6428  # K31
$ws.Cells.Item(31, 13).Value = -1295.6428  # M31
$ws.Cells.Item(34, 8).Value = 1590.6428  # H34
$ws.Cells.Item(34, 9).Value = 1590.6428  # I34
$ws.Cells.Item(34, 11).Value = 1590.6428  # K34
$ws.Cells.Item(34, 13).Value = -1388.6428  # M34
$ws.Cells.Item(54, 8).Value = 152092  # H54
$ws.Cells.Item(54, 10).Value = 152092  # J54
$ws.Cells.Item(54, 12).Value = 152092  # L54
$ws.Cells.Item(54, 14).Value = -153408  # N54
$ws.Cells.Item(56, 8).Value = 16499.5  # H56
$ws.Cells.Item(56, 9).Value = 7999  # I56
$ws.Cells.Item(56, 11).Value = 7999  # K56
$ws.Cells.Item(56, 13).Value = -7154  # M56
$ws.Cells.Item(59, 8).Value = 27820.572  # H59
$ws.Cells.Item(59, 9).Value = 16501.334  # I59
$ws.Cells.Item(59, 10).Value = 30907.637  # J59
$ws.Cells.Item(59, 11).Value = 16501.334  # K59
$ws.Cells.Item(59, 12).Value = 30907.637  # L59
$ws.Cells.Item(59, 13).Value = -15356.334  # M59
$ws.Cells.Item(59, 14).Value = -33197.637  # N59
$ws.Cells.Item(60, 8).Value = 20676.666  # H60
$ws.Cells.Item(60, 9).Value = 21218  # I60
$ws.Cells.Item(60, 10).Value = 20000  # J60
$ws.Cells.Item(60, 11).Value = 21218  # K60
$ws.Cells.Item(60, 12).Value = 20000  # L60
$ws.Cells.Item(60, 13).Value = -20707  # M60
$ws.Cells.Item(60, 14).Value = -21022  # N60
$ws.Cells.Item(68, 8).Value = 30000  # H68
$ws.Cells.Item(68, 9).Value = 0  # I68
$ws.Cells.Item(68, 11).Value = 0  # K68
$ws.Cells.Item(71, 8).Value = 30000  # H71
$ws.Cells.Item(71, 9).Value = 0  # I71
$ws.Cells.Item(71, 11).Value = 0  # K71
$ws.Cells.Item(68, 13).ClearContents()  # M68
$ws.Cells.Item(71, 13).ClearContents()  # M71

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 500113  # H11
$ws.Cells.Item(11, 10).Value = 150.66667  # J11
$ws.Cells.Item(11, 12).Value = 452.00001  # L11
$ws.Cells.Item(11, 14).Value = -732.00001  # N11

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7499.5  # H70
$ws.Cells.Item(70, 9).Value = 7499.5  # I70
$ws.Cells.Item(70, 11).Value = 7499.5  # K70
$ws.Cells.Item(70, 13).Value = -7229.5  # M70
$ws.Cells.Item(73, 8).Value = 7499.5  # H73
$ws.Cells.Item(73, 9).Value = 7499.5  # I73
$ws.Cells.Item(73, 11).Value = 7499.5  # K73
$ws.Cells.Item(73, 13).Value = -6563.5  # M73
$ws.Cells.Item(122, 8).Value = 1549  # H122
$ws.Cells.Item(122, 9).Value = 1549  # I122
$ws.Cells.Item(122, 11).Value = 4647  # K122
$ws.Cells.Item(122, 13).Value = -2197  # M122

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4226.909  # H132
$ws.Cells.Item(132, 9).Value = 3249.5  # I132
$ws.Cells.Item(132, 11).Value = 9748.5  # K132
$ws.Cells.Item(132, 13).Value = -7218.5  # M132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1814.5294  # H107
$ws.Cells.Item(107, 9).Value = 1986.5555  # I107
$ws.Cells.Item(107, 10).Value = 1621  # J107
$ws.Cells.Item(107, 11).Value = 5959.666499999999  # K107
$ws.Cells.Item(107, 12).Value = 4863  # L107
$ws.Cells.Item(107, 13).Value = -4039.666499999999  # M107
$ws.Cells.Item(107, 14).Value = -8703  # N107
$ws.Cells.Item(136, 8).Value = 3055.6  # H136
$ws.Cells.Item(136, 9).Value = 3822.4285  # I136
$ws.Cells.Item(136, 11).Value = 11467.2855  # K136
$ws.Cells.Item(136, 13).Value = -8917.2855  # M136
